$wb = $excel.ActiveWorkbook

# YDS sheet - append new weeks play-by-play yardage log data
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '6 3 0 13 3 6 5 -1 16 3 0 8 8 1 -2 8 5 1 7 8 10 2 3 0 6 20 2 5 3 2 5 -1 4 7 1 1 3 8 10 1 3 3 1 2 0 4 4 3 2 2 5 6 11 10 -1 12 23 0 1 0 6 17 0 5 9 6 3 -5 6 3 -7 45 3 -3 4 4 -1 10 9 8 11 2 2 4 3 2 3 2 0 3 -3 1 2 -1 3 6 0 2 3 2 2 4 3 3 1 4 7 3 1 1 2 3 1 3 4 -3 3 -1 7 1 0 -1 5 11 0 0 1 12 0 9 2 7 1 3 1 3 1 4 3 1 4 1 3 5 11 2 6 3 0 3 5 3 9 -1 5 1 4 5 38 1 4 1 7 1 -1 3 5 1 -2 0 0 4 1 2 -1 5 2 2 6 0 -2 -3 4 3 3 5 5 1 4 5 5 3 2 4 6 2 6 0 0 0 -5 5 -1 -2 3 2 57 -1 5 3 3 11 9 2 0 6 -3 4 4 10 4 8 4 8 -1 3 12 13 8 13 4 1 2 4 10 5 7 3 0 6 4 1 13 2 2 2 -1 3 9 -3 3 3 80 3 -1 0 3 3 1 4 13 9 11 5 5 12 0 2 5 -2 3 5 9 -2 1 2 19 3 8 2 16 8 10 9 0 1 3 7 3 2 1 2 6 3 4 2 7 6 -1 3 11 14 2 8 3 5 1 1 7 3 9 3 2 -1 4 0 1 2 -2 2 26 1 1 -3 2 2 4 4 3 -1 8 6 12 0 3 12 6 3 3 0 23 3 5 0 2 0 0 3 3 3 2 6 -1 6 8 2 4 2 1 5 -1 3 5 6 4 3 -1 2 7 0 4 5 3 0 8 13 2 3 3 0 2 4 6 4 1 5 3 0 0 3 2 0 1 7 41 0 14 2 -3 -2 1 3 5 10 3 3 0 9 3 1 6 4 4 6 1 3 8 5 4 3 2 2 3 7 8 6 11 14 5 0 -1 0 6 6 5 11 7 2 0 -3 6 1 1 4 4 -1 1 10 6 4 16 1 7 0 -2 1 6 4 2 1 5 4 1 9 4 1 4 4 10 1 4 -1 5 1 11 9 1 4 10 2 3 2 10 1 4 3 3 24 4 -1 5 7 2 23 8 -3 4 -1 4 6 1 1 3 0 5 11 14 0 0 6 8 9 3 3 9 7 0 4 8 8 1 0 0 6 2 0 7 1 6 6 5 0 1 0 1 1 9 2 0 0 1 25 5 4 1 3 7 3 4 2 13 4 2 2 11 3 6 1 2 14 3 16 0 29 1 1 6 4 0 6 13 5 12 6 0 6 0 8 8 1 2 -1 4 -1 -1 1 1 8 2 6 16 -1 4 22 9 -2 1 3 6 9 6 3 2 3 4 16 3 9 3 16 2 0 15 -3 -4 4 -4 4 6 3 0 14 8 4 0 0 2 3 8 21 2 8 1 6 2 8 1 0 6 2 2 3 2 10 3 1 12 5 2 15 3 8 8 9 2 11 5 -6 3 10 7 5 6 4 5 5 -1 6 11 2 0 7 2 0 0 0 1 2 4 8 4 0 2 6 0 0 3 3 5 4 5 12 -1 9 1 0 1 0 4 9 2 8 7 3 1 13 0 24 -1 7 -2 1 2 2 3 3 8 1 0 5 5 4 6 5 8 2 1 5 1 20 1 1 3 8 3 8 7 0 20 4 5 2 2 3 11 17 1 9 7 13 -1 1 9 3 4 3 3 -2 16 -1 8 0 4 4 11 -1 2 0 1 -7 5 -2 23 3 4 11 2 2 -2 8 3 2 -2 0 0 1 -2 1 12 8 3 -4 4 4 -1 9 2 20 3 0 -1 1'
$wsYDS.Range("C2").Value = '2 -1 19 2 13 2 0 1 1 1 4 4 5 3 21 10 6 1 2 -2 14 14 2 3 2 1 7 -1 5 18 0 6 3 6 3 2 7 1 -1 3 2 12 6 0 3 2 5 7 1 4 7 35 11 1 6 2 16 10 1 2 4 3 16 3 4 1 1 3 2 5 12 -2 3 -1 12 9 3 0 2 1 -1 3 6 -3 0 5 3 -3 3 8 16 1 -4 12 0 4 1 4 -2 1 4 4 11 4 2 6 0 5 3 3 6 37 10 -1 2 -2 -1 3 10 11 3 4 3 1 3 1 12 10 10 -1 8 0 1 6 1 4 1 8 0 1 3 1 11 11 5 6 4 9 9 0 -3 16 7 -1 13 6 3 5 12 1 25 10 11 1 9 -3 5 1 -4 3 -2 5 -1 2 -3 2 3 9 3 1 20 2 -1 12 3 -2 3 1 2 3 1 4 3 -3 6 1 2 3 11 10 2 20 3 1 5 2 5 0 3 0 6 2 2 4 0 3 5 3 1 0 5 6 1 26 0 0 2 2 1 1 5 4 5 2 0 3 2 1 5 0 3 -1 6 2 0 0 2 3 1 0 1 8 4 13 6 0 7 0 14 4 1 3 1 14 1 9 8 1 4 0 3 5 1 10 2 1 7 2 2 5 7 3 11 4 8 2 0 17 4 4 1 13 2 16 1 -4 1 4 10 2 0 4 15 4 2 1 3 5 2 5 3 3 2 3 0 1 3 2 5 9 0 4 2 10 -1 3 5 4 7 14 3 5 6 5 1 4 1 3 3 2 3 2 9 3 1 0 10 0 13 7 0 6 6 2 7 20 4 10 11 8 4 18 1 11 2 1 2 2 9 1 3 0 3 3 -2 25 10 0 8 2 0 4 11 1 2 2 3 4 -1 10 -1 4 1 2 4 4 3 4 4 4 9 3 4 1 4 6 14 8 4 2 5 2 6 10 -1 9 4 7 -2 3 6 0 1 1 2 -1 1 4 5 4 0 -2 3 0 2 25 8 1 5 1 0 6 1 2 7 0 -2 0 5 1 5 8 9 5 2 2 15 3 3 8 2 3 2 3 1 8 -1 -3 4 9 5 1 5 10 3 3 5 2 4 3 2 12 1 4 2 2 10 12 4 1 3 13 -1 4 1 11 6 0 3 0 12 1 4 7 29 14 3 5 3 19 1 6 2 -6 4 6 5 4 2 5 5 8 4 5 3 -1 1 3 11 2 -2 2 4 9 2 -1 2 4 6 6 3 1 3 -4 0 3 11 3 8 6 5 0 5 -1 5 5 4 2 7 1 1 0 1 1 8 16 6 0 3 3 0 5 8 2 6 36 1 -6 28 1 5 1 5 1 0 7 6 8 3 9 3 15 2 4 5 12 3 9 4 11 6 8 12 2 -1 1 12 19 12 9 5 3 1 2 -3 7 1 1 27 1 -1 3 0 4 2 3 11 5 27 1 5 39 -1 12 -2 0 1 1 3 12 10 2 2 7 2 0 5 8 4 5 4 -2 11 0 11 2 0 0 3 1 2 1 13 -1 3 1 5 10 1 -1 -3 4 1 4 2 5 4 2 19 0 3 8 2 3 4 5 2 5 4 6 5 0 2 1 6 -1 2 2 3 4 -4 0 0 2 4 0 3 7 1 8 4 2 4 2 6 11 19 3 4 1 8 5 0 9 8 21 2 7 1 2 -4 5 10 7 2 14 -1 -6 10 2 5 4 9 0 -6 9 3 3 6 3 0 -4 2 15 1 2 3 6 11 0 11 5 8 11 1 3 14 0 5 3 6 6 2 1 9 7 1 4 7 10 2 3 10 9 2 0 -1 8 0 5 19 11 -2 6 2 6 0 3 1 3 3 5 3 2 2 4 2 1 3 3 -1 25 5 0 2 3 2 6 5 7 3 10 6 2 32 3 7 32 3 3 7 2'
$wsYDS.Range("B3").Value = '13 12 10 19 8 11 12 22 19 3 0 18 22 3 5 2 7 18 13 16 1 27 17 15 2 28 18 9 -1 14 16 12 2 12 15 4 8 5 9 8 11 5 7 5 19 13 12 17 16 2 8 5 11 20 3 7 7 30 -8 4 6 3 4 37 6 2 28 8 0 6 5 6 33 27 0 17 5 5 11 11 3 12 4 5 4 8 4 7 25 4 2 10 22 16 10 5 24 13 8 7 9 -6 4 12 13 25 7 8 12 3 9 8 9 8 6 8 5 2 10 -1 17 4 1 9 17 11 2 7 4 2 8 12 11 7 10 13 39 6 13 6 8 8 18 6 7 1 23 7 8 8 7 2 38 7 2 3 8 14 6 18 9 14 6 7 42 3 4 7 5 7 19 6 5 0 8 4 9 11 7 50 24 -4 2 9 5 5 16 5 6 7 9 28 11 -3 3 2 6 2 15 17 8 24 4 4 5 4 14 9 7 6 34 5 5 7 22 18 3 -4 23 4 5 5 8 5 27 13 1 6 4 16 1 2 12 11 8 0 7 7 3 5 24 9 8 7 21 3 -1 7 10 -4 18 3 8 20 16 7 12 12 9 6 13 1 13 10 12 12 6 16 8 4 7 10 9 4 9 3 16 17 21 10 13 10 7 10 15 19 9 13 -1 6 4 21 6 4 11 7 9 16 1 14 5 4 14 12 7 9 7 7 10 5 4 6 9 11 12 8 22 16 16 3 2 11 35 18 13 10 7 4 24 8 14 12 10 5 16 9 16 13 22 7 35 3 21 12 5 30 12 13 6 12 9 12 10 12 24 2 8 7 8 13 22 2 8 6 6 7 8 7 5 5 1 6 6 1 0 5 14 10 5 2 11 2 53 3 5 5 14 2 6 4 12 4 6 5 3 9 28 6 1 11 11 2 21 15 3 8 18 5 5 18 7 6 6 19 4 9 5 10 4 19 -2 10 9 7 3 10 9 8 6 7 10 8 3 10 1 11 5 11 9 9 10 5 6 10 17 -5 3 3 11 2 36 10 5 11 10 13 21 5 0 11 9 9 18 11 10 4 4 21 64 8 8 28 32 6 27 7 4 18 7 2 2 15 8 12 8 9 7 13 14 15 11 20 3 5 3 12 8 11 20 21 8 10 5 15 7 6 5 5 5 3 6 3 17 4 12 4 10 22 20 6 22 12 -6 11 8 5 14 6 11 11 16 6 8 10 4 2 -6 19 22 6 2 9 10 20 5 3 6 25 8 21 17 50 10 11 20 28 5 9 39 16 11 12 6 29 22 22 23 60 10 2 3 12 9 4 23 3 11 13 8 7 12 33 2 17 52 17 16 8 9 0 10 19 11 16 12 7 13 12 7 12 3 5 6 3 10 7 6 10 34 11 2 1 7 12 11 6 6 2 7 19 0 11 19 8 3 6 46 6 54 6 2 5 6 13 5 11 4 3 22 5 12 27 6 8 9 0 11 6 2 9 11 8 12 11 23 11 6 13 3 16 10 10 17 19 4 1 10 0 8 15 9 15 7 5 16 2 15 16 10 16 14 11 5 30 10 14 15'
$wsYDS.Range("C3").Value = '20 26 18 21 7 7 16 24 15 14 16 4 9 4 9 4 5 11 5 3 32 5 17 15 24 12 13 -3 7 7 11 20 12 6 4 2 12 10 11 8 22 4 9 4 63 1 16 2 9 14 6 12 5 2 13 2 24 13 22 3 12 13 6 11 36 13 15 11 13 12 11 8 3 -6 13 27 16 -1 8 11 17 2 35 9 2 0 9 4 9 10 1 8 6 14 1 11 7 23 9 16 12 4 22 12 5 5 16 1 20 9 38 2 13 19 -4 13 11 -3 16 14 7 7 4 18 9 3 34 2 6 -5 13 21 6 10 3 11 8 12 9 7 -1 1 1 47 14 4 5 5 11 5 12 4 8 16 5 12 11 13 3 9 20 4 -3 10 9 10 6 14 38 6 17 9 40 6 12 30 2 -2 2 8 9 12 12 6 17 23 -3 4 9 9 13 11 14 4 2 8 54 7 22 14 6 8 21 11 2 8 5 12 8 15 11 5 8 8 15 4 2 5 0 39 24 7 21 19 5 4 15 14 14 3 25 49 5 22 15 12 7 9 24 9 23 6 5 14 22 17 5 13 25 8 13 6 5 11 6 3 48 5 8 3 3 10 8 4 15 12 4 3 15 11 9 8 21 23 3 9 5 9 2 3 16 26 30 13 3 9 14 5 9 20 6 12 7 11 7 9 0 18 5 2 4 5 14 10 7 20 3 3 6 6 8 1 14 3 19 8 34 9 3 15 9 3 15 3 7 72 11 9 13 15 6 7 9 8 3 17 14 6 5 13 1 4 9 0 17 11 6 17 7 24 0 6 38 10 7 3 11 6 5 10 5 5 14 12 10 9 7 67 1 18 4 17 17 1 6 19 56 6 37 3 6 17 15 9 2 2 7 12 10 3 12 1 2 14 8 7 14 14 22 5 10 15 42 7 13 13 21 7 0 3 18 17 7 18 17 6 4 9 13 26 16 23 15 7 21 5 6 33 17 2 20 11 22 5 9 4 11 9 24 13 25 17 3 8 9 9 3 12 8 29 9 29 4 4 9 6 17 5 6 18 3 2 4 6 3 12 7 10 13 4 0 6 32 1 12 5 11 9 4 10 12 14 13 41 8 -4 6 4 21 4 11 9 10 46 2 9 14 0 9 8 9 7 12 0 34 21 -4 19 11 6 23 16 6 15 50 3 83 15 7 11 19 23 2 26 9 11 4 0 7 3 1 3 42 10 9 4 10 2 11 4 12 22 13 11 1 16 5 9 8 -2 11 11 2 3 12 9 15 2 9 8 -2 3 9 4 12 7 21 6 29 0 3 6 16 39 9 10 5 4 -6 9 2 0 -1 17 8 0 11 17 15 7 10 12 23 8 7 12 13 0 4 14 6 6 32 0 7 12 25 12 20 6 13 6 18 38 12 16 8 23 16 12 6 6 8 11 6 4 7 3 5 3 0 2 18 12 4 9 7 7 7 13 7 41 0 9 8 9 13 13 26 3 3 24 9 -3 4 15'

# OFF sheet - updated season totals
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 375
$wsOFF.Range("D2").Value = 27
$wsOFF.Range("E2").Value = 26
$wsOFF.Range("F2").Value = 135
$wsOFF.Range("G2").Value = 98
$wsOFF.Range("J2").Value = 58
$wsOFF.Range("N2").Value = 46
$wsOFF.Range("O2").Value = 56
$wsOFF.Range("P2").Value = 29
$wsOFF.Range("B3").Value = 14
$wsOFF.Range("C3").Value = 354
$wsOFF.Range("E3").Value = 72
$wsOFF.Range("F3").Value = 180
$wsOFF.Range("G3").Value = 81
$wsOFF.Range("H3").Value = 62
$wsOFF.Range("I3").Value = 116
$wsOFF.Range("J3").Value = 101
$wsOFF.Range("L3").Value = 600
$wsOFF.Range("M3").Value = 376
$wsOFF.Range("Q3").Value = 1082

# DEF sheet - updated season totals
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 399
$wsDEF.Range("D2").Value = 20
$wsDEF.Range("E2").Value = 21
$wsDEF.Range("F2").Value = 135
$wsDEF.Range("G2").Value = 119
$wsDEF.Range("N2").Value = 37
$wsDEF.Range("O2").Value = 47
$wsDEF.Range("C3").Value = 320
$wsDEF.Range("E3").Value = 56
$wsDEF.Range("F3").Value = 184
$wsDEF.Range("G3").Value = 64
$wsDEF.Range("H3").Value = 55
$wsDEF.Range("I3").Value = 107
$wsDEF.Range("J3").Value = 124
$wsDEF.Range("L3").Value = 543
$wsDEF.Range("M3").Value = 351
$wsDEF.Range("Q3").Value = 1049

# ST sheet - updated totals and appended per-kick logs
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 144
$wsST.Range("D2").Value = 127
$wsST.Range("F2").Value = 229
$wsST.Range("G2").Value = 221
$wsST.Range("H2").Value = 7
$wsST.Range("I2").Value = 4
$wsST.Range("J2").Value = 113
$wsST.Range("K2").Value = 106
$wsST.Range("B3").Value = 74
$wsST.Range("B4").Value = '68 68 66 69 65 62 67 66 63 63 63 65 63 65 62 61 47 62 40 64 54 38 61 47 57 47 55 63 61 60 63 65 62 66 65 62 63 65 64 64 60 70 66 62 67 66 66 62 62 66 64 65 61 56 59 69 62 67 66 66 60 63 64 60 61 66 62 66 46 54'
$wsST.Range("B5").Value = '31 35 22 34 14 12 32 24 25 21 19 35 23 18 28 18 14 31 10 23 14 8 16 10 23 13 9 15 10 20 19 30 18 23 30 6 28 20 21 15 23 27 18 24 23 21 26 15 15 20 23 17 29 21 19 25 25 35 27 22 21 32 18 19 18 33 23 40 13 16'
$wsST.Range("B6").Value = '29 36 45 23 30 30 10 26 25 40 0 26 8 21 28 32 47 40 28 22 21 23 17 8 32 4 20 7 26 4 7 31 16 21 20 17 45 12 5 5 16 8 18 20 30 22 0 20 33 26 42 50 24 13 19 10 22 39 23 31 33 24 32 13 26 27 30 27 31 10 34 24 30 27 25 16 19 26 32 17 17 0 25 19 15 22 27 40 0 42 22 14 18 34 27 22 30 20 8 14'
$wsST.Range("D3").Value = '45 54 33 55 42 38 50 51 33 49 54 48 64 18 50 33 49 42 62 59 41 38 47 43 59 35 39 49 54 43 55 47 44 49 53 45 47 55 54 64 53 36 51 40 54 50 31 39 47 47 41 44 41 32 53 37 48 28 46 40 44 58 40 35 53 49 42 47 46 31 25 45 42 36 42 40 47 40 61 61 52 50 63 48 62 55 44 60 43 53 52 44 42 48 43 43 27 50 28 38 44 41 44 43 49 54 31 53 40 34 50 28 46 51 39 53 38 40 49 42 54 44 45 72 64 28 43'
$wsST.Range("D4").Value = '0 0 0 16 0 0 11 0 0 7 11 0 14 0 3 0 4 14 0 21 0 0 12 0 8 0 0 9 0 0 15 0 0 -1 0 42 11 14 16 21 12 0 0 0 0 0 0 0 0 0 0 11 9 0 0 0 0 0 7 0 0 0 0 0 7 7 0 4 7 0 0 11 0 0 0 0 15 0 14 13 11 23 9 14 24 7 0 22 0 0 0 0 9 16 0 0 0 43 0 0 0 21 7 14 6 3 0 7 12 0 0 0 9 0 0 18 0 0 13 0 15 11 3 2 0 0 15'
$wsST.Range("D5").Value = '0 0 34 0 97 0 6 -1 0 0 0 0 4 28 0 0'

# TURNS sheet - updated turnover totals
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("D3").Value = 15
$wsTURNS.Range("E3").Value = 16

# PEN sheet - updated penalty totals
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 39
$wsPEN.Range("B3").Value = 22

